# Auto-generated edit script: update crypto price/volume data
# per commit "Updated cryptos list on Fri May 17 15:00:01 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.472.06"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "3.085.03"
$ws.Range("E3").Value = "  +4.61%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "'580.55"
$ws.Range("E5").Value = "  +1.91%  "
$ws.Range("D6").Value = "'167.74"
$ws.Range("E6").Value = "  +5.55%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "3.079.57"
$ws.Range("E8").Value = "  +4.67%  "
$ws.Range("D9").Value = "'0.522"
$ws.Range("E9").Value = "  +1.75%  "
$ws.Range("D10").Value = "'6.68"
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("E11").Value = "  +2.01%  "
$ws.Range("E12").Value = "  +6.46%  "
$ws.Range("D13").Value = "'0.0000249"
$ws.Range("E13").Value = "  +2.55%  "
$ws.Range("D14").Value = "'36.70"
$ws.Range("E14").Value = "  +8.21%  "
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").Value = "3.596.14"
$ws.Range("E16").Value = "  +4.52%  "
$ws.Range("D17").Value = "66.408.71"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "'7.17"
$ws.Range("E18").Value = "  +4.79%  "
$ws.Range("D19").Value = "3.083.62"
$ws.Range("E19").Value = "  +4.26%  "
$ws.Range("D20").Value = "'15.99"
$ws.Range("E20").Value = "  +17.60%  "
$ws.Range("D21").Value = "'464.31"
$ws.Range("E21").Value = "  +4.52%  "
$ws.Range("E22").Value = "  +6.47%  "
$ws.Range("D23").Value = "'7.45"
$ws.Range("E23").Value = "  +5.21%  "
$ws.Range("D24").Value = "'83.17"
$ws.Range("E24").Value = "  +1.92%  "
$ws.Range("D25").Value = "'12.78"
$ws.Range("E25").Value = "  +5.92%  "
$ws.Range("D26").Value = "'2.28"
$ws.Range("E26").Value = "  +4.03%  "
$ws.Range("E27").Value = "  +1.69%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").Value = "'8.06"
$ws.Range("E29").Value = "  +1.09%  "
$ws.Range("D30").Value = "'2.41"
$ws.Range("E30").Value = "  +1.95%  "
$ws.Range("E31").Value = "  +3.94%  "
$ws.Range("E32").Value = "  +5.46%  "
$ws.Range("D33").Value = "'28.50"
$ws.Range("E33").Value = "  +6.16%  "
$ws.Range("E34").Value = "  +6.40%  "
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  +2.75%  "
$ws.Range("D37").Value = "'5.89"
$ws.Range("E37").Value = "  +4.35%  "
$ws.Range("D38").Value = "'48.29"
$ws.Range("E38").Value = "  +12.97%  "
$ws.Range("D39").Value = "'50.06"
$ws.Range("E39").Value = "  +1.80%  "
$ws.Range("E40").Value = "  +3.92%  "
$ws.Range("E41").Value = "  +5.09%  "
$ws.Range("E42").Value = "  +3.35%  "
$ws.Range("D43").Value = "'2.89"
$ws.Range("E43").Value = "  +4.14%  "
$ws.Range("D44").Value = "'8.64"
$ws.Range("E44").Value = "  +4.26%  "
$ws.Range("E45").Value = "  +2.86%  "
$ws.Range("D46").Value = "'385.43"
$ws.Range("E46").Value = "  +1.77%  "
$ws.Range("D47").Value = "2.770.62"
$ws.Range("E47").Value = "  +2.76%  "
$ws.Range("D48").Value = "'134.67"
$ws.Range("E48").Value = "  +3.03%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "'24.54"
$ws.Range("E50").Value = "  +7.48%  "
$ws.Range("E51").Value = "  +5.78%  "
